$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "WARNING: not found -> $old"
    }
}

Replace-Text 'Unveiling the Quantum Realm: A Path to Technological Singularity' 'Mathematics: A Journey Through Numbers and Logic'
Replace-Text 'Dr' 'Prof'
Replace-Text ' Amelia Jacobson' ' Marcus Adams'
Replace-Text 'amelia' 'marcus'
Replace-Text 'jacobson@quantuminstitute' 'adams@eduinstitute'
Replace-Text 'In the realm of science, the quest to unravel the mysteries of the quantum realm holds immense significance, promising transformative potential across diverse fields' 'Mathematics, the language of the universe, holds immense power to unveil the secrets of the cosmos'
Replace-Text ' Quantum mechanics, with its perplexing phenomena like superposition, entanglement, and wave-particle duality, presents a new frontier of exploration, beckoning us to transcend the limitations of classical physics' ' It is a subject that captivates the minds of young and old alike, enriching our understanding of nature and enabling us to unravel the mysteries it holds'
Replace-Text ' As we delve deeper into this enigmatic domain, we stand at the cusp of technological breakthroughs that could herald an era of unparalleled progress' ' From the intricate patterns of fractals to the elegance of calculus, mathematics invites us on an exploration of the fundamental principles governing our world'
Replace-Text 'The intricate dance of subatomic particles offers a glimpse into a world governed by probabilities and uncertainties, where quantum properties defy our conventional notions of reality' 'In the realm of mathematics, we investigate the complex interplay of numbers, symbols, and equations'
Replace-Text ' Harnessing these quantum effects promises exponential leaps in computing power, enabling us to solve previously intractable problems and revolutionize fields such as cryptography, optimization, and artificial intelligence' ' We unlock the mysteries of geometry, traversing through the landscapes of shapes and angles'
Replace-Text ' Beyond computing, quantum technologies are poised to revolutionize materials science, medicine, and communication, opening up new avenues for innovation and societal transformation' ' We discover the beauty of algebra, manipulating expressions and equations to unveil hidden relationships. And as we delve into calculus, we uncover the intricate dance of change, exploring the rates and patterns that define the universe''s dynamic processes'
Replace-Text 'While the path to technological singularity, a hypothetical point where technological advancements reach an irreversible and self-sustaining crescendo, remains shrouded in uncertainty, the exploration of the quantum realm offers a compelling path forward' 'Mathematics is not merely about abstract concepts; it is a tool that empowers us to solve real-world problems'
Replace-Text ' This journey, though fraught with challenges, beckons us to transcend the boundaries of our current understanding and forge a new era of scientific discovery and technological marvels' ' From predicting the trajectory of a rocket to analyzing financial data, mathematics equips us with the skills to navigate the complexities of our world. It fosters critical thinking, analytical reasoning, and problem-solving abilities, preparing us for success in various fields and endeavors'
Replace-Text 'Our exploration of the quantum realm, with its enigmatic phenomena and profound implications, holds the key to unlocking a future of boundless possibilities' 'Mathematics stands as a powerful and versatile tool that unveils the mysteries of the universe and equips us with essential skills for navigating the complexities of life'
Replace-Text ' From transformative computing technologies to groundbreaking advancements in materials science and medicine, the uncharted territory of quantum mechanics beckons us to transcend the limitations of classical physics' ' It invites us on an exploration of numbers, symbols, and equations, captivating our minds with its intricate patterns and elegant structures'
Replace-Text ' While the path to technological singularity remains uncertain, the pursuit of quantum knowledge offers a compelling roadmap to a world of unprecedented technological prowess, ushering in an era of innovation that will redefine the very fabric of our existence' ' As we delve deeper into mathematics, we unlock the secrets of geometry, algebra, and calculus, gaining a profound understanding of the fundamental principles that govern our world. Mathematics is more than a subject; it is a language that empowers us to solve real-world problems and make informed decisions, shaping our lives and contributing to the advancement of society'

# Append a trailing empty paragraph at the end of the document body
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

